# Auto-generated script to update horarios workbook with newly scraped schedule data
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:50:23"
$ws.Cells.Item(3, 1).Value = "Total filas: 95"
$newRows = @(
  @{R=55; A="07:17:57"; B="07:31"; C="11_ETCHEVERRY"; D=14; E="LP1912"},
  @{R=56; A="07:17:57"; B="07:31"; C="84_COLONIA URQUIZA-ESC 49"; D=14; E="LP1912"},
  @{R=57; A="07:17:57"; B="07:31"; C="16_SANTA ANA"; D=14; E="LP1912"},
  @{R=65; A="07:50:23"; B="07:51"; C="215D_EL PATO"; D=1; E="LP1912"},
  @{R=66; A="07:50:23"; B="07:59"; C="23_HERNANDEZ"; D=9; E="LP1912"},
  @{R=68; A="07:50:23"; B="08:07"; C="11_ETCHEVERRY"; D=17; E="LP1912"},
  @{R=69; A="07:50:23"; B="08:10"; C="16_SANTA ANA"; D=20; E="LP1912"},
  @{R=70; A="07:17:57"; B="08:11"; C="15_ABASTO"; D=54; E="LP1912"},
  @{R=71; A="07:50:23"; B="08:12"; C="15_ABASTO"; D=22; E="LP1912"},
  @{R=72; A="07:50:23"; B="08:13"; C="10_OLMOS"; D=23; E="LP1912"},
  @{R=73; A="07:17:57"; B="08:20"; C="26_HERNANDEZ"; D=63; E="LP1912"},
  @{R=74; A="07:50:23"; B="08:21"; C="26_HERNANDEZ"; D=31; E="LP1912"},
  @{R=75; A="07:50:23"; B="08:22"; C="16_P MOR-SANTA ANA"; D=32; E="LP1912"},
  @{R=76; A="07:17:57"; B="08:22"; C="215B_EL PATO"; D=65; E="LP1912"},
  @{R=77; A="07:50:23"; B="08:23"; C="215B_EL PATO"; D=33; E="LP1912"},
  @{R=78; A="07:17:57"; B="08:26"; C="84_COLONIA URQUIZA-ESC 49"; D=69; E="LP1912"},
  @{R=79; A="07:50:23"; B="08:27"; C="84_COLONIA URQUIZA-ESC 49"; D=37; E="LP1912"},
  @{R=80; A="07:50:23"; B="08:34"; C="23_HERNANDEZ"; D=44; E="LP1912"},
  @{R=81; A="06:52:31"; B="08:35"; C="23_HERNANDEZ"; D=103; E="LP1912"},
  @{R=82; A="07:17:57"; B="08:41"; C="81_EL PELIGRO"; D=84; E="LP1912"},
  @{R=83; A="07:50:23"; B="08:42"; C="81_EL PELIGRO"; D=52; E="LP1912"},
  @{R=84; A="07:50:23"; B="08:43"; C="14_ABASTO"; D=53; E="LP1912"},
  @{R=85; A="07:17:57"; B="08:53"; C="17_ROMERO"; D=96; E="LP1912"},
  @{R=86; A="07:50:23"; B="08:53"; C="10_OLMOS"; D=63; E="LP1912"},
  @{R=87; A="07:50:23"; B="08:54"; C="17_ROMERO"; D=64; E="LP1912"},
  @{R=88; A="07:50:23"; B="09:01"; C="215A_EL PATO"; D=71; E="LP1912"},
  @{R=89; A="07:50:23"; B="09:03"; C="11_ETCHEVERRY"; D=73; E="LP1912"},
  @{R=90; A="07:50:23"; B="09:10"; C="16_P MOR-SANTA ANA"; D=80; E="LP1912"},
  @{R=91; A="07:17:57"; B="09:16"; C="27_EL RETIRO"; D=119; E="LP1912"},
  @{R=92; A="07:50:23"; B="09:17"; C="27_EL RETIRO"; D=87; E="LP1912"},
  @{R=93; A="07:50:23"; B="09:21"; C="26_HERNANDEZ"; D=91; E="LP1912"},
  @{R=94; A="07:50:23"; B="09:23"; C="11_ETCHEVERRY"; D=93; E="LP1912"},
  @{R=95; A="07:50:23"; B="09:23"; C="17_ROMERO"; D=93; E="LP1912"},
  @{R=96; A="07:50:23"; B="09:31"; C="16_SANTA ANA"; D=101; E="LP1912"},
  @{R=97; A="07:50:23"; B="09:32"; C="15_ABASTO"; D=102; E="LP1912"},
  @{R=98; A="07:50:23"; B="09:33"; C="10_OLMOS"; D=103; E="LP1912"},
  @{R=99; A="07:50:23"; B="09:42"; C="215C_EL PATO"; D=112; E="LP1912"},
  @{R=100; A="07:50:23"; B="09:43"; C="14_ABASTO"; D=113; E="LP1912"}
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.R, 1).Value = $row.A
  $ws.Cells.Item($row.R, 2).Value = $row.B
  $ws.Cells.Item($row.R, 3).Value = $row.C
  $ws.Cells.Item($row.R, 4).Value = $row.D
  $ws.Cells.Item($row.R, 5).Value = $row.E
}

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:50:23"
$ws.Cells.Item(3, 1).Value = "Total filas: 20"
$newRows = @(
  @{R=21; A="07:50:23"; B="07:51"; C="215D_EL PATO"; D=1; E="LP1912"},
  @{R=23; A="07:50:23"; B="08:23"; C="215B_EL PATO"; D=33; E="LP1912"},
  @{R=24; A="07:50:23"; B="09:01"; C="215A_EL PATO"; D=71; E="LP1912"},
  @{R=25; A="07:50:23"; B="09:42"; C="215C_EL PATO"; D=112; E="LP1912"}
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.R, 1).Value = $row.A
  $ws.Cells.Item($row.R, 2).Value = $row.B
  $ws.Cells.Item($row.R, 3).Value = $row.C
  $ws.Cells.Item($row.R, 4).Value = $row.D
  $ws.Cells.Item($row.R, 5).Value = $row.E
}

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:50:23"
$ws.Cells.Item(3, 1).Value = "Total filas: 22"
$newRows = @(
  @{R=19; A="07:50:23"; B="07:53"; C="215A_LA PLATA"; D=3; E="L6173"},
  @{R=20; A="06:52:31"; B="08:07"; C="215C_LA PLATA"; D=75; E="L6203"},
  @{R=21; A="07:17:57"; B="08:09"; C="215C_LA PLATA"; D=52; E="L6203"},
  @{R=22; A="07:50:23"; B="08:14"; C="215C_LA PLATA"; D=24; E="L6203"},
  @{R=23; A="06:52:31"; B="08:30"; C="215A_LA PLATA"; D=98; E="L6173"},
  @{R=24; A="07:17:57"; B="08:34"; C="215A_LA PLATA"; D=77; E="L6173"},
  @{R=25; A="07:50:23"; B="08:35"; C="215A_LA PLATA"; D=45; E="L6173"},
  @{R=26; A="07:17:57"; B="09:08"; C="215D_LA PLATA"; D=111; E="L6203"},
  @{R=27; A="07:50:23"; B="09:09"; C="215D_LA PLATA"; D=79; E="L6203"}
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.R, 1).Value = $row.A
  $ws.Cells.Item($row.R, 2).Value = $row.B
  $ws.Cells.Item($row.R, 3).Value = $row.C
  $ws.Cells.Item($row.R, 4).Value = $row.D
  $ws.Cells.Item($row.R, 5).Value = $row.E
}
